# Generate Report for Handback
#
# The "bec9db9b-4814-4cc0-a64f-6ea3ce7bf266" file (row 3 in every sheet) has
# been handed back and is now in sync with en-US, for both the zh-cn and
# de-de locales. Update the Status columns accordingly and stamp the
# "Latest Handback DateTime" for each locale-specific sheet.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn (col B) / de-de (col C) status for row 3 ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $statusHandedBack
$overview.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet: Status (col C) + Latest Handback DateTime (col H) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusHandedBack
$zhcn.Range("H3").Value = "2016-03-11 08:28:59"

# --- de-de sheet: Status (col C) + Latest Handback DateTime (col H) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusHandedBack
$dede.Range("H3").Value = "2016-03-11 08:29:05"
